$d = $word.ActiveDocument

# The requested change trims the document-wide style defaults
# (w:docDefaults in word/styles.xml) down to just the handful of
# properties that aren't already implied by the OOXML spec defaults.
# There's no dedicated Word object-model surface for "docDefaults", so
# we round-trip the whole package through WordOpenXML and rewrite that
# one element directly.

$oldDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/><w:shd w:val="clear" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:pBdr><w:top w:val="nil" w:sz="0" w:space="0"/><w:left w:val="nil" w:sz="0" w:space="0"/><w:bottom w:val="nil" w:sz="0" w:space="0"/><w:right w:val="nil" w:sz="0" w:space="0"/><w:between w:val="nil" w:sz="0" w:space="0"/></w:pBdr><w:shd w:val="clear" w:fill="auto"/><w:spacing w:before="0" w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="0" w:right="0" w:firstLine="0"/><w:contextualSpacing w:val="0"/><w:jc w:val="left"/></w:pPr></w:pPrDefault></w:docDefaults>'

$newDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

$xml = $d.WordOpenXML

if ($xml.Contains($oldDocDefaults)) {
    $xml = $xml.Replace($oldDocDefaults, $newDocDefaults)
} else {
    # Fallback: locate the docDefaults element irrespective of exact
    # attribute ordering/formatting and replace just that span.
    $startTag = "<w:docDefaults>"
    $endTag = "</w:docDefaults>"
    $start = $xml.IndexOf($startTag)
    $end = $xml.IndexOf($endTag)
    if ($start -ge 0 -and $end -ge 0) {
        $end = $end + $endTag.Length
        $before = $xml.Substring(0, $start)
        $after = $xml.Substring($end)
        $xml = $before + $newDocDefaults + $after
    }
}

$d.WordOpenXML = $xml

Write-Output "docDefaults rewritten"
